# Generate Report for Handback
# Updates the handoff/handback timestamps for the file that was just
# handed back (c03eaaa9-6ef0-4fde-8007-3eed47ecb17f.md), on the Overview
# sheet and on each per-locale report sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for c03eaaa9...md (row 3)
$overview.Range("G3").Value = "2016-09-07 01:06:11"

# zh-cn report sheet, row 3 (c03eaaa9...md):
#   Correspond Handoff Datetime
#   Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-09-07 01:06:00"
$zhcn.Range("K3").Value = "2016-09-07 01:06:29"

# de-de report sheet, row 3 (c03eaaa9...md):
#   Correspond Handoff Datetime
#   Correspond Handback DateTime
$dede.Range("H3").Value = "2016-09-07 01:06:11"
$dede.Range("K3").Value = "2016-09-07 01:06:37"
